$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Ford F-150 XL" -> "2019 Ford F-150 XL" (rest of the row is unchanged)
$ws.Range("A4").Value = "2019 Ford F-150 XL"

# New rows 5-8 with additional car data
$ws.Range("A5").Value = "2019 Audi S4 Sedan"
$ws.Range("B5").Value = 51195
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = "sedan"
$ws.Range("E5").Value = "white"
$ws.Range("F5").Value = "gas"
$ws.Range("G5").Value = 349
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = "all"
$ws.Range("J5").Value = "audi"

$ws.Range("A6").Value = "2019 Audi Q3 SUV"
$ws.Range("B6").Value = 34700
$ws.Range("C6").Value = 27
$ws.Range("D6").Value = "suv"
$ws.Range("E6").Value = "black"
$ws.Range("F6").Value = "gas"
$ws.Range("G6").Value = 228
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = "all"
$ws.Range("J6").Value = "audi"

$ws.Range("A7").Value = "2019 Hyundai Sonata"
$ws.Range("B7").Value = 22650
$ws.Range("C7").Value = 37
$ws.Range("D7").Value = "sedan"
$ws.Range("E7").Value = "blue"
$ws.Range("F7").Value = "gas"
$ws.Range("G7").Value = 185
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = "front"
$ws.Range("J7").Value = "hyundai"

$ws.Range("A8").Value = "2019 Hyundai IONIQ Hybrid"
$ws.Range("B8").Value = 20650
$ws.Range("C8").Value = 59
$ws.Range("D8").Value = "sedan"
$ws.Range("E8").Value = "blue"
$ws.Range("F8").Value = "hybrid"
$ws.Range("G8").Value = 139
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = "front"
$ws.Range("J8").Value = "hyundai"

# Update selection to match the state after entering the new data
$ws.Range("A9:F9").Select()
